$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.241.68'
$ws.Range("E2").Value = '  -3.37%  '

$ws.Range("D3").Value = '3.367.85'
$ws.Range("E3").Value = '  -4.10%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.73'
$ws.Range("E5").Value = '  -3.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '124.71'
$ws.Range("E6").Value = '  -7.35%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.367.53'
$ws.Range("E8").Value = '  -4.08%  '

$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.22'
$ws.Range("E10").Value = '  -5.53%  '

$ws.Range("E11").Value = '  -4.41%  '

$ws.Range("E12").Value = '  -3.80%  '

$ws.Range("D13").Value = '3.940.90'
$ws.Range("E13").Value = '  -4.08%  '

$ws.Range("E14").Value = '  -1.00%  '

$ws.Range("D15").Value = '3.367.28'
$ws.Range("E15").Value = '  -4.10%  '

$ws.Range("E16").Value = '  -6.14%  '

$ws.Range("D17").Value = '62.311.45'
$ws.Range("E17").Value = '  -3.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.29'
$ws.Range("E18").Value = '  -5.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.19'
$ws.Range("E19").Value = '  -8.12%  '

$ws.Range("E20").Value = '  -2.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.04'
$ws.Range("E21").Value = '  -4.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '371.41'
$ws.Range("E22").Value = '  -6.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("E23").Value = '  -4.55%  '

$ws.Range("D24").Value = '3.503.46'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.22'

$ws.Range("E27").Value = '  -10.75%  '

$ws.Range("E28").Value = '  -0.08%  '

$ws.Range("E29").Value = '  -7.15%  '

$ws.Range("E30").Value = '  -7.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.76'
$ws.Range("E31").Value = '  -5.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  -6.77%  '

$ws.Range("D34").Value = '3.399.68'
$ws.Range("E34").Value = '  -3.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.148'
$ws.Range("E35").Value = '  -6.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.60'
$ws.Range("E36").Value = '  -3.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.17'
$ws.Range("E37").Value = '  -3.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.68'
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.61'
$ws.Range("E39").Value = '  -5.19%  '

$ws.Range("E40").Value = '  -5.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0748'
$ws.Range("E41").Value = '  -5.09%  '

$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.762'
$ws.Range("E43").Value = '  -6.13%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.36'
$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("E45").Value = '  -5.29%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.53'
$ws.Range("E46").Value = '  -9.26%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.53'
$ws.Range("E47").Value = '  -7.64%  '

$ws.Range("E48").Value = '  -9.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.58'
$ws.Range("E49").Value = '  -3.43%  '

$ws.Range("D50").Value = '2.236.46'
$ws.Range("E50").Value = '  -5.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.834'
$ws.Range("E51").Value = '  -8.14%  '
